$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A149").Value = 111734024
$ws.Range("B149").Value = 88950
$ws.Range("D149").Value = 'NT'
$ws.Range("E149").Value = 256756
$ws.Range("F149").Value = 'Blek fingersvamp'
$ws.Range("G149").Value = 'Ramaria pallida'
$ws.Range("H149").Value = '(Schaeff.) Ricken'
$ws.Range("I149").Value = '2'
$ws.Range("J149").Value = 'mycel'
$ws.Range("AC149").Value = '13+8 fruktkroppar'
$ws.Range("A151").Value = 111734063
$ws.Range("B151").Value = 88946
$ws.Range("D151").Value = 'VU'
$ws.Range("E151").Value = 256335
$ws.Range("F151").Value = 'Taggfingersvamp'
$ws.Range("G151").Value = 'Ramaria karstenii'
$ws.Range("H151").Value = '(Sacc. & P.Syd.) Corner'
$ws.Range("I151").Value = '5'
$ws.Range("J151").Value = 'fruktkroppar'
$ws.Range("Q151").Value = 537734.3589083791
$ws.Range("R151").Value = 6656722.215160147
$ws.Range("AJ151").Value = 'gran'
$ws.Range("AK151").Value = 'Picea abies'
$ws.Range("AO151").Value = 'Picea abies'
$ws.Range("A152").Value = 111733878
$ws.Range("B152").Value = 90678
$ws.Range("D152").Value = 'LC'
$ws.Range("E152").Value = 4366
$ws.Range("F152").Value = 'Skarp dropptaggsvamp'
$ws.Range("G152").Value = 'Hydnellum peckii'
$ws.Range("H152").Value = 'Banker'
$ws.Range("I152").ClearContents()
$ws.Range("J152").ClearContents()
$ws.Range("Q152").Value = 537749.1720232533
$ws.Range("R152").Value = 6656839.262154824
$ws.Range("AC152").ClearContents()
$ws.Range("A153").Value = 111733658
$ws.Range("B153").Value = 90678
$ws.Range("D153").Value = 'LC'
$ws.Range("E153").Value = 4366
$ws.Range("F153").Value = 'Skarp dropptaggsvamp'
$ws.Range("G153").Value = 'Hydnellum peckii'
$ws.Range("H153").Value = 'Banker'
$ws.Range("I153").ClearContents()
$ws.Range("J153").ClearContents()
$ws.Range("Q153").Value = 537747.6347874232
$ws.Range("R153").Value = 6657038.059664771
$ws.Range("A154").Value = 111734358
$ws.Range("B154").Value = 88966
$ws.Range("D154").Value = 'NT'
$ws.Range("E154").Value = 5754
$ws.Range("F154").Value = 'Gultoppig fingersvamp'
$ws.Range("G154").Value = 'Ramaria testaceoflava'
$ws.Range("H154").Value = '(Bres.) Corner'
$ws.Range("I154").Value = '17'
$ws.Range("Q154").Value = 537695.459656042
$ws.Range("R154").Value = 6656709.327821301
$ws.Range("A156").Value = 111734405
$ws.Range("B156").Value = 90018
$ws.Range("D156").Value = 'LC'
$ws.Range("E156").Value = 1339
$ws.Range("F156").Value = 'Brandticka'
$ws.Range("G156").Value = 'Pycnoporellus fulgens'
$ws.Range("H156").Value = '(Fr.) Donk'
$ws.Range("I156").Value = '2'
$ws.Range("Q156").Value = 537687.4847450438
$ws.Range("R156").Value = 6656706.248840296
$ws.Range("AJ156").Value = 'gran'
$ws.Range("AK156").Value = 'Picea abies'
$ws.Range("AO156").Value = 'högstubbe # Picea abies'
$ws.Range("A158").Value = 111734390
$ws.Range("B158").Value = 89665
$ws.Range("E158").Value = 73
$ws.Range("F158").Value = 'Veckticka'
$ws.Range("G158").Value = 'Flavidoporia pulvinascens'
$ws.Range("H158").Value = '(Pilát) Audet'
$ws.Range("I158").ClearContents()
$ws.Range("J158").ClearContents()
$ws.Range("Q158").Value = 537687.4847450438
$ws.Range("R158").Value = 6656706.248840296
$ws.Range("AJ158").Value = 'asp'
$ws.Range("AK158").Value = 'Populus tremula'
$ws.Range("AO158").Value = 'grov asplåga # Populus tremula'
$ws.Range("A159").Value = 111733654
$ws.Range("B159").Value = 90687
$ws.Range("E159").Value = 5964
$ws.Range("F159").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G159").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H159").Value = '(L.:Fr.) P.Karst.'
$ws.Range("A160").Value = 111734100
$ws.Range("B160").Value = 88967
$ws.Range("D160").Value = 'DD'
$ws.Range("E160").Value = 6039940
$ws.Range("F160").Value = 'Mandarinfingersvamp'
$ws.Range("G160").Value = 'Ramaria tridentina'
$ws.Range("H160").Value = 'Schild'
$ws.Range("I160").Value = '2'
$ws.Range("Q160").Value = 537710.2246313525
$ws.Range("R160").Value = 6656733.457844303
$ws.Range("AJ160").Value = 'gran'
$ws.Range("AK160").Value = 'Picea abies'
$ws.Range("AO160").Value = 'Picea abies'
$ws.Range("A161").Value = 111733686
$ws.Range("I161").Value = '13'
$ws.Range("J161").Value = 'fruktkroppar'
$ws.Range("Q161").Value = 537776.5022700967
$ws.Range("R161").Value = 6656955.434839563
$ws.Range("A162").Value = 111733729
$ws.Range("B162").Value = 90662
$ws.Range("D162").Value = 'LC'
$ws.Range("E162").Value = 4363
$ws.Range("F162").Value = 'Zontaggsvamp'
$ws.Range("G162").Value = 'Hydnellum concrescens'
$ws.Range("H162").Value = '(Pers.) Banker'
$ws.Range("I162").Value = '10'
$ws.Range("Q162").Value = 537784.8865190516
$ws.Range("R162").Value = 6656918.555543642
$ws.Range("A163").Value = 111734445
$ws.Range("B163").Value = 85210
$ws.Range("D163").Value = 'LC'
$ws.Range("E163").Value = 3624
$ws.Range("F163").Value = 'Strimspindling'
$ws.Range("G163").Value = 'Cortinarius glaucopus'
$ws.Range("H163").Value = '(Schaeff. : Fr.) Fr.'
$ws.Range("I163").Value = '20'
$ws.Range("Q163").Value = 537682.2670869593
$ws.Range("R163").Value = 6656678.219876539
$ws.Range("A164").Value = 111733781
$ws.Range("B164").Value = 88946
$ws.Range("D164").Value = 'VU'
$ws.Range("E164").Value = 256335
$ws.Range("F164").Value = 'Taggfingersvamp'
$ws.Range("G164").Value = 'Ramaria karstenii'
$ws.Range("H164").Value = '(Sacc. & P.Syd.) Corner'
$ws.Range("I164").Value = '25'
$ws.Range("Q164").Value = 537777.7813424434
$ws.Range("R164").Value = 6656879.518825633
$ws.Range("A165").Value = 111734039
$ws.Range("B165").Value = 88967
$ws.Range("D165").Value = 'DD'
$ws.Range("E165").Value = 6039940
$ws.Range("F165").Value = 'Mandarinfingersvamp'
$ws.Range("G165").Value = 'Ramaria tridentina'
$ws.Range("H165").Value = 'Schild'
$ws.Range("I165").Value = '3'
$ws.Range("J165").Value = 'fruktkroppar'
$ws.Range("Q165").Value = 537725.9133506906
$ws.Range("R165").Value = 6656765.090555903
$ws.Range("AC165").Value = 'små'
$ws.Range("A166").Value = 111733999
$ws.Range("I166").Value = '12'
$ws.Range("Q166").Value = 537735.408403003
$ws.Range("R166").Value = 6656815.142909602
$ws.Range("A167").Value = 111734434
$ws.Range("B167").Value = 88909
$ws.Range("D167").Value = 'VU'
$ws.Range("E167").Value = 720
$ws.Range("F167").Value = 'Violgubbe'
$ws.Range("G167").Value = 'Gomphus clavatus'
$ws.Range("H167").Value = '(Pers.) Gray'
$ws.Range("I167").Value = '12'
$ws.Range("Q167").Value = 537683.5943855111
$ws.Range("R167").Value = 6656695.218654346
$ws.Range("AC167").Value = 'ring ca 1,5 m i diameter'
$ws.Range("A168").Value = 111734658
$ws.Range("B168").Value = 90151
$ws.Range("E168").Value = 366
$ws.Range("F168").Value = 'Kandelabersvamp'
$ws.Range("G168").Value = 'Artomyces pyxidatus'
$ws.Range("H168").Value = '(Pers.) Jülich'
$ws.Range("I168").Value = '9'
$ws.Range("Q168").Value = 537685.1971427263
$ws.Range("R168").Value = 6656734.200801319
$ws.Range("AJ168").Value = 'asp'
$ws.Range("AK168").Value = 'Populus tremula'
$ws.Range("AO168").Value = 'grov asplåga # Populus tremula'
$ws.Range("A169").Value = 111734629
$ws.Range("B169").Value = 88966
$ws.Range("E169").Value = 5754
$ws.Range("F169").Value = 'Gultoppig fingersvamp'
$ws.Range("G169").Value = 'Ramaria testaceoflava'
$ws.Range("H169").Value = '(Bres.) Corner'
$ws.Range("I169").Value = '1'
$ws.Range("J169").Value = 'fruktkroppar'
$ws.Range("Q169").Value = 537663.3456201598
$ws.Range("R169").Value = 6656669.03414992
$ws.Range("AC169").ClearContents()
$ws.Range("A170").Value = 111734315
$ws.Range("B170").Value = 90151
$ws.Range("D170").Value = 'NT'
$ws.Range("E170").Value = 366
$ws.Range("F170").Value = 'Kandelabersvamp'
$ws.Range("G170").Value = 'Artomyces pyxidatus'
$ws.Range("H170").Value = '(Pers.) Jülich'
$ws.Range("I170").Value = '20'
$ws.Range("Q170").Value = 537701.3922571414
$ws.Range("R170").Value = 6656716.382399381
$ws.Range("AJ170").Value = 'asp'
$ws.Range("AK170").Value = 'Populus tremula'
$ws.Range("AO170").Value = 'grov murken asplåga # Populus tremula'
